$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the weekly observation data which gets reshuffled across rows.
$cols = @(4, 12, 13, 14, 15, 16, 17, 19, 20)

# Mapping of destination row -> source row (data in the new row comes from the
# old content of the source row), as derived from the commit's diff.
$rowMap = @{
    2  = 14
    3  = 2
    4  = 10
    5  = 4
    6  = 13
    7  = 12
    8  = 8
    9  = 11
    10 = 15
    11 = 9
    12 = 7
    13 = 5
    14 = 6
    15 = 3
}

# Snapshot all current values (rows 2-15, relevant columns) before writing
# anything, since several rows form permutation cycles and would otherwise
# clobber data that is still needed.
$snapshot = @{}
for ($r = 2; $r -le 15; $r++) {
    foreach ($c in $cols) {
        $snapshot["$r,$c"] = $ws.Cells.Item($r, $c).Value()
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value = $snapshot["$srcRow,$c"]
    }
}
